# Generate Report for Handback
#
# This script updates the localization-status report to reflect that the
# zh-cn and de-de handbacks have completed and are now in sync with en-US:
#   - The "Ready for handoff" status becomes "Handed back: in sync with en-US"
#     (this text is shared across the Overview sheet and the per-locale sheets).
#   - The zh-cn and de-de sheets get their "Latest Target File" / "Latest
#     Handback File" columns populated with the handed-back xliff file names,
#     and a hyperlink is added on the target-file cell.
#   - The zh-cn and de-de "Latest Handback DateTime" is stamped.
#   - A few columns are widened to better fit the longer values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$mdFileName = "d8737fb6-aa9c-452c-afa9-e35f63e21a94.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/a10dda9853ec3673a5a4976da853810100a85485/e2e/d8737fb6-aa9c-452c-afa9-e35f63e21a94.md"

$handedBackStatus = "Handed back: in sync with en-US"

# --- Status: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# (shared by Overview!E2, Overview!F2, zh-cn!C2, de-de!C2)
$wsOverview.Range("E2").Value = $handedBackStatus
$wsOverview.Range("F2").Value = $handedBackStatus
$wsZhCn.Range("C2").Value = $handedBackStatus
$wsDeDe.Range("C2").Value = $handedBackStatus

# --- zh-cn: Latest Target File / Latest Handback File / Handback DateTime ---
$wsZhCn.Range("I2").Value = $mdFileName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null

$wsZhCn.Range("J2").Value = "d8737fb6-aa9c-452c-afa9-e35f63e21a94.07e02ea9d1f1514fa21061e889d9b7927b342ac4.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-01 15:24:18"

# --- de-de: Latest Target File / Latest Handback File / Handback DateTime ---
$wsDeDe.Range("I2").Value = $mdFileName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdFileName) | Out-Null

$wsDeDe.Range("J2").Value = "d8737fb6-aa9c-452c-afa9-e35f63e21a94.07e02ea9d1f1514fa21061e889d9b7927b342ac4.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-01 15:24:26"

# --- Column width adjustments ---
# Overview: Status columns (zh-cn / de-de) widened to fit the new, longer text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# zh-cn / de-de: Status column (C) widened, and Latest Target File / Latest
# Handback File columns (I / J) widened to fit the file names.
$wsZhCn.Columns.Item(3).ColumnWidth = 29.15
$wsZhCn.Columns.Item(9).ColumnWidth = 39.15
$wsZhCn.Columns.Item(10).ColumnWidth = 39.15

$wsDeDe.Columns.Item(3).ColumnWidth = 29.15
$wsDeDe.Columns.Item(9).ColumnWidth = 39.15
$wsDeDe.Columns.Item(10).ColumnWidth = 39.15
